# Database Design Document.xlsx edit
# Commit message: "changes in ddd file for changes in datatype instead of
# datetime change to timestamp"
#
# Summary of the edit:
#  1. Fix a typo on the "Scope" sheet: "verious" -> "various".
#  2. Rename the "RecordCreateDateTime" / "RecordModifiedDateTime" columns
#     (and their "DATETIME" datatype cells) to "RecordCreateTIMESTAMP" /
#     "RecordModifiedTIMESTAMP" with datatype "TIMESTAMP" throughout the
#     "Tables Description" and "Authentication Tables" sheets.
#  3. The active sheet moves from "Scope" to "Tables Description", and the
#     selected cells on a few sheets change.

$wb = $excel.ActiveWorkbook

$wsScope = $wb.Worksheets.Item("Scope")
$wsTablesDesc = $wb.Worksheets.Item("Tables Description")
$wsAuth = $wb.Worksheets.Item("Authentication Tables")

# --- 1. Typo fix on the "Scope" sheet -------------------------------------
$wsScope.Range("F7").Value = "To identify various tables needed for database"

# --- 2. DATETIME -> TIMESTAMP renames --------------------------------------

# "Tables Description" sheet: F-column holds the field name,
# H-column holds the corresponding datatype.
$tdCreateRefs = @("F9","F20","F36","F54","F60","F74","F85","F98","F104","F115","F120","F138","F151","F165","F178")
foreach ($ref in $tdCreateRefs) {
    $wsTablesDesc.Range($ref).Value = "RecordCreateTIMESTAMP"
}

$tdModifiedRefs = @("F10","F21","F37","F55","F75","F86","F99","F105","F139","F152","F166","F179")
foreach ($ref in $tdModifiedRefs) {
    $wsTablesDesc.Range($ref).Value = "RecordModifiedTIMESTAMP"
}

$tdDatatypeRefs = @("H9","H10","H20","H21","H36","H37","H54","H55","H60","H74","H75","H85","H86","H98","H99","H104","H105","H115","H120","H138","H139","H151","H152","H165","H166","H171","H178","H179")
foreach ($ref in $tdDatatypeRefs) {
    $wsTablesDesc.Range($ref).Value = "TIMESTAMP"
}

# "Authentication Tables" sheet
$authCreateRefs = @("F6","F12","F18")
foreach ($ref in $authCreateRefs) {
    $wsAuth.Range($ref).Value = "RecordCreateTIMESTAMP"
}

$authModifiedRefs = @("F7","F19")
foreach ($ref in $authModifiedRefs) {
    $wsAuth.Range($ref).Value = "RecordModifiedTIMESTAMP"
}

$authDatatypeRefs = @("H6","H7","H12","H18","H19")
foreach ($ref in $authDatatypeRefs) {
    $wsAuth.Range($ref).Value = "TIMESTAMP"
}

# --- 3. Selection / active sheet updates -----------------------------------
$wsScope.Range("F8").Select() | Out-Null
$wsAuth.Range("F3").Select() | Out-Null

$wsTablesDesc.Activate() | Out-Null
$wsTablesDesc.Range("H10").Select() | Out-Null
